$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.045.32"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.590.84"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.33"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.75"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  +13.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.414"
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("E9").Value = "  +6.83%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.588.32"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.87"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.255.69"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.815.36"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.591.34"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.72"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.08"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.530"
$ws.Range("E22").Value = "  +10.28%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.43"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "509.48"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("E26").Value = "  +5.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "97.09"
$ws.Range("E27").Value = "  +5.92%  "
$ws.Range("E28").Value = "  +5.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.784.08"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  +10.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.56"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.55"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.93"
$ws.Range("E37").Value = "  +11.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "624.78"
$ws.Range("E38").Value = "  +10.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.572"
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  +12.35%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.916"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.85"
$ws.Range("E44").Value = "  +6.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.82"
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.55"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.37"
$ws.Range("E49").Value = "  -8.21%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.32"
$ws.Range("E51").Value = "  +3.99%  "
